$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.335162643758053
$ws.Range("C2").Value = 0.300772252573978
$ws.Range("D2").Value = 0.2211918794791075
$ws.Range("E2").Value = 0.1743160135658499
$ws.Range("F2").Value = 1.186764554642728
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.1819088135301783
$ws.Range("O2").Value = 2.688228375106689
$ws.Range("B3").Value = 1.193101117243259
$ws.Range("C3").Value = 0.2627201377795529
$ws.Range("D3").Value = 0.2147181916055416
$ws.Range("E3").Value = 0.1699470556310203
$ws.Range("F3").Value = 1.191627199630915
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.1779213376851487
$ws.Range("O3").Value = 2.715148437091955
$ws.Range("B4").Value = 1.105774502175848
$ws.Range("C4").Value = 0.2392675052076072
$ws.Range("D4").Value = 0.2108020149816667
$ws.Range("E4").Value = 0.1673444571717972
$ws.Range("F4").Value = 1.195574216601983
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.1755854790393556
$ws.Range("O4").Value = 2.73416583233103
$ws.Range("B5").Value = 1.070164852877838
$ws.Range("C5").Value = 0.2296886941643379
$ws.Range("D5").Value = 0.20922101095708
$ws.Range("E5").Value = 0.1663040130096469
$ws.Range("F5").Value = 1.197423874964493
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.1746618415502112
$ws.Range("O5").Value = 2.742539517129671
$ws.Range("B6").Value = 1.064250543345622
$ws.Range("C6").Value = 0.2280968480349088
$ws.Range("D6").Value = 0.2089593878207694
$ws.Range("E6").Value = 0.166132465087756
$ws.Range("F6").Value = 1.197745563055022
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.1745101770415047
$ws.Range("O6").Value = 2.743967593513531
$ws.Range("B7").Value = 1.105294350018198
$ws.Range("C7").Value = 0.2391384089929716
$ws.Range("D7").Value = 0.2107806326531261
$ws.Range("E7").Value = 0.1673303438086258
$ws.Range("F7").Value = 1.195598185755138
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.1755729082310964
$ws.Range("O7").Value = 2.734276238828414
$ws.Range("B8").Value = 1.28620170744523
$ws.Range("C8").Value = 0.2876705791924223
$ws.Range("D8").Value = 0.2189476337364766
$ws.Range("E8").Value = 0.1727930170370797
$ws.Range("F8").Value = 1.188241355749241
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.1805105650978476
$ws.Range("O8").Value = 2.696992880599851
$ws.Range("B9").Value = 1.640101407823238
$ws.Range("C9").Value = 0.3821194928075329
$ws.Range("D9").Value = 0.2354249270776023
$ws.Range("E9").Value = 0.1841392161275692
$ws.Range("F9").Value = 1.181467067778968
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.1910881103853939
$ws.Range("O9").Value = 2.643705853848985
$ws.Range("B10").Value = 1.89952970033039
$ws.Range("C10").Value = 0.4510502535528076
$ws.Range("D10").Value = 0.2478085709550868
$ws.Range("E10").Value = 0.1928621402399386
$ws.Range("F10").Value = 1.181191195775611
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.1994094103585979
$ws.Range("O10").Value = 2.616756535286981
$ws.Range("B11").Value = 2.017413814856639
$ws.Range("C11").Value = 0.4823048200645985
$ws.Range("D11").Value = 0.2535018369542144
$ws.Range("E11").Value = 0.1969146005760578
$ws.Range("F11").Value = 1.1820943959369
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.2033154431349828
$ws.Range("O11").Value = 2.60717058789507
$ws.Range("B12").Value = 2.062033205660384
$ws.Range("C12").Value = 0.4941249350540602
$ws.Range("D12").Value = 0.2556662529640619
$ws.Range("E12").Value = 0.1984612832275943
$ws.Range("F12").Value = 1.182584944859912
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.2048119639956951
$ws.Range("O12").Value = 2.603927060630582
$ws.Range("B13").Value = 2.052424575137138
$ws.Range("C13").Value = 0.4915799520407518
$ws.Range("D13").Value = 0.2551997309854812
$ws.Range("E13").Value = 0.1981276395424914
$ws.Range("F13").Value = 1.182472681038618
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.2044888872728166
$ws.Range("O13").Value = 2.604608391963723
$ws.Range("B14").Value = 2.021085105990437
$ws.Range("C14").Value = 0.4832775782346062
$ws.Range("D14").Value = 0.2536797350491469
$ws.Range("E14").Value = 0.1970416046603205
$ws.Range("F14").Value = 1.182131773554346
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.2034382140099211
$ws.Range("O14").Value = 2.606895984060571
$ws.Range("B15").Value = 2.001886015826756
$ws.Range("C15").Value = 0.4781901232918244
$ws.Range("D15").Value = 0.2527497972755555
$ws.Range("E15").Value = 0.1963779524327691
$ws.Range("F15").Value = 1.181942318267559
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.2027969124447964
$ws.Range("O15").Value = 2.608347589271801
$ws.Range("B16").Value = 1.891822893805568
$ws.Range("C16").Value = 0.4490055834085069
$ws.Range("D16").Value = 0.247437698095041
$ws.Range("E16").Value = 0.1925989976332332
$ws.Range("F16").Value = 1.181152920933698
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.1991565731298834
$ws.Range("O16").Value = 2.617436969079591
$ws.Range("B17").Value = 1.824267782389597
$ws.Range("C17").Value = 0.4310751692004828
$ws.Range("D17").Value = 0.2441941553564391
$ws.Range("E17").Value = 0.1903023182619847
$ws.Range("F17").Value = 1.180932553014216
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.1969542709747287
$ws.Range("O17").Value = 2.623699182089439
$ws.Range("B18").Value = 1.785399640029425
$ws.Range("C18").Value = 0.4207524670915745
$ws.Range("D18").Value = 0.242334198550509
$ws.Range("E18").Value = 0.1889892708140266
$ws.Range("F18").Value = 1.180902592998322
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.1956989167814811
$ws.Range("O18").Value = 2.627552556476019
$ws.Range("B19").Value = 1.772237513225036
$ws.Range("C19").Value = 0.4172557464395368
$ws.Range("D19").Value = 0.2417054217744408
$ws.Range("E19").Value = 0.1885460603415652
$ws.Range("F19").Value = 1.18090905446131
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.1952758237506487
$ws.Range("O19").Value = 2.628900379027783
$ws.Range("B20").Value = 1.831460418920585
$ws.Range("C20").Value = 0.4329848904790197
$ws.Range("D20").Value = 0.2445388527821137
$ws.Range("E20").Value = 0.1905459819535906
$ws.Range("F20").Value = 1.18094598981348
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.1971875346156509
$ws.Range("O20").Value = 2.623006513256058
$ws.Range("B21").Value = 2.030290846391722
$ws.Range("C21").Value = 0.4857166068559877
$ws.Range("D21").Value = 0.2541259649037073
$ws.Range("E21").Value = 0.1973602711460956
$ws.Range("F21").Value = 1.182227870585621
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.2037463498735832
$ws.Range("O21").Value = 2.606213558173465
$ws.Range("B22").Value = 2.160116125128013
$ws.Range("C22").Value = 0.5200903207177134
$ws.Range("D22").Value = 0.2604411889329867
$ws.Range("E22").Value = 0.201884341647812
$ws.Range("F22").Value = 1.183931659366522
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.2081343037889809
$ws.Range("O22").Value = 2.597491932432945
$ws.Range("B23").Value = 2.090837707934838
$ws.Range("C23").Value = 0.501752805087051
$ws.Range("D23").Value = 0.2570661430362975
$ws.Range("E23").Value = 0.1994633144937055
$ws.Range("F23").Value = 1.182942877077252
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.2057830783651013
$ws.Range("O23").Value = 2.601939959950641
$ws.Range("B24").Value = 1.828208720656505
$ws.Range("C24").Value = 0.432121549894589
$ws.Range("D24").Value = 0.2443830001204503
$ws.Range("E24").Value = 0.1904357987219001
$ws.Range("F24").Value = 1.180939613760842
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.1970820425541149
$ws.Range("O24").Value = 2.623318880439513
$ws.Range("B25").Value = 1.544460254482601
$ws.Range("C25").Value = 0.3566480535240544
$ws.Range("D25").Value = 0.2309183347159376
$ws.Range("E25").Value = 0.1810018623999667
$ws.Range("F25").Value = 1.182476715971831
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.1881303467702224
$ws.Range("O25").Value = 2.655987105482751
